# Weekly crypto-symbol refresh (GitHub Actions bot): updated prices and the
# "Hora" (hour) column for row 2..51, plus a couple of "Bestin24h"/"Worstin24h"
# badge moves between the CEJI (row 41) and KickToken (row 42) volume labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is a cell address -> new text value. The sheet stores these as
# plain text (prices/hours are inline strings, not numbers), so we write the
# value with a leading apostrophe to force text entry (mirrors typing '274.64
# into Excel) and then clear the resulting number-format flag so the cell
# ends up with no explicit style, exactly like the other untouched text cells.
$updates = [ordered]@{
    "D2" = "274.64"
    "G2" = "21"
    "D3" = "21.12"
    "G3" = "21"
    "D4" = "6.245"
    "G4" = "21"
    "D5" = "0.06206"
    "G5" = "21"
    "D6" = "3.553"
    "G6" = "21"
    "D7" = "1.522"
    "G7" = "21"
    "D8" = "6.549"
    "G8" = "21"
    "D9" = "0.8241"
    "G9" = "21"
    "D10" = "0.1647"
    "G10" = "21"
    "D11" = "0.08268"
    "G11" = "21"
    "D12" = "0.03489"
    "G12" = "21"
    "D13" = "0.03124"
    "G13" = "21"
    "D14" = "0.09161"
    "G14" = "21"
    "D15" = "3.759"
    "G15" = "21"
    "D16" = "0.001625"
    "G16" = "21"
    "D17" = "0.04664"
    "G17" = "21"
    "D18" = "0.006236"
    "G18" = "21"
    "D19" = "0.006201"
    "G19" = "21"
    "D20" = "0.001066"
    "G20" = "21"
    "D21" = "0.0001496"
    "G21" = "21"
    "D22" = "3.724"
    "G22" = "21"
    "G23" = "21"
    "D24" = "0.01395"
    "G24" = "21"
    "G25" = "21"
    "D26" = "0.1248"
    "G26" = "21"
    "G27" = "21"
    "D28" = "0.0002730"
    "G28" = "21"
    "G29" = "21"
    "G30" = "21"
    "G31" = "21"
    "G32" = "21"
    "G33" = "21"
    "G34" = "21"
    "G35" = "21"
    "G36" = "21"
    "G37" = "21"
    "G38" = "21"
    "G39" = "21"
    "D40" = "0.04730"
    "G40" = "21"
    "D41" = "0.005287"
    "E41" = "40CEJICEJI"
    "G41" = "21"
    "D42" = "0.007032"
    "E42" = "41KickTokenKICKBestin24h"
    "G42" = "21"
    "D43" = "0.1117"
    "G43" = "21"
    "D44" = "0.01137"
    "G44" = "21"
    "D45" = "0.00006043"
    "G45" = "21"
    "D46" = "0.00000000748"
    "G46" = "21"
    "D47" = "0.7213"
    "G47" = "21"
    "D48" = "0.001388"
    "G48" = "21"
    "D49" = "0.00001895"
    "G49" = "21"
    "D50" = "0.01237"
    "G50" = "21"
    "G51" = "21"
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.ClearFormats()
}
